$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct the floating point precision of the existing last row (row 21, column A)
$ws.Range("A21").Value = 45866.91692689815

# Append the new row of sensor data (row 22)
$ws.Range("A22").Value = 45866.95866324242
$ws.Range("B22").Value = 2025
$ws.Range("C22").Value = 31
$ws.Range("D22").Value = 12.89
$ws.Range("E22").Value = 89.76000000000001
$ws.Range("F22").Value = 0
$ws.Range("G22").Value = 0.99
$ws.Range("H22").Value = "ENE"
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = "23:00:28"

# Match the date/time number formatting + style used by the rest of column A
$ws.Range("A22").NumberFormat = $ws.Range("A21").NumberFormat
